$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.185631155967712
$ws.Range("B1").Value = 2.347034692764282
$ws.Range("C1").Value = 3.951726198196411
$ws.Range("D1").Value = 3.009632110595703
$ws.Range("E1").Value = 1.134935855865479
